# TFEC_renewables.xlsx update:
#  - Residential biomass end-use rows (2-6) move from the "Renewables" group
#    to the "Fossil fuels" group (category C).
#  - A new "Solar" row (RES_CWH_SOLAR) is inserted after the existing Solar
#    block, pushing every subsequent row down by one.
#  - The old "RES_CWH_EL_001" residential datafile id is renamed/replaced by
#    "RES_CWH_KER_001" (row stays in the Kerosene / Fossil fuels bucket).
#  - Active selection moves to C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, make room for the new row by duplicating the formatting of the last
# populated row (C35) onto the new last row (C36), then write the full,
# final table (rows 1-36) so every value/grouping matches the target state.
$ws.Range("A35:C35").Copy()
$ws.Range("A36:C36").PasteSpecial(-4122)

# Prime the two brand-new datafile ids in the same order the author typed
# them (Kerosene rename before the new Solar row), so they land in the
# shared-string table in that order too.
$ws.Cells.Item(15, 2).Value = "RES_CWH_KER_001"
$ws.Cells.Item(13, 2).Value = "RES_CWH_SOLAR"

$target = @(
    @(1, "VISUALIZATION", "OSEMOSYS", "Renewable"),
    @(2, "Biomass", "RES_CWH_BIO_001", "Fossil fuels"),
    @(3, "Biomass", "COM_CWH_BIO_001", "Fossil fuels"),
    @(4, "Biomass", "TRA_AN_BIO_001", "Fossil fuels"),
    @(5, "Biomass", "OTH_ALL_BIO_001", "Fossil fuels"),
    @(6, "Biomass", "IND_ALL_BIO_001", "Fossil fuels"),
    @(7, "Biomass", "PWR_BIO_001", "Renewables"),
    @(8, "Geothermal", "PWR_GEO_001", "Renewables"),
    @(9, "Hydro", "PWR_HYD_001", "Renewables"),
    @(10, "Solar", "PWR_SOL_E001", "Renewables"),
    @(11, "Solar", "PWR_SOL_002", "Renewables"),
    @(12, "Solar", "PWR_SOL_003", "Renewables"),
    @(13, "Solar", "RES_CWH_SOLAR", "Renewables"),
    @(14, "Wind", "PWR_WND_E001", "Renewables"),
    @(15, "Kerosene", "RES_CWH_KER_001", "Fossil fuels"),
    @(16, "Natural Gas", "RES_CWH_NGS_001", "Fossil fuels"),
    @(17, "Natural Gas", "COM_CWH_NGS_001", "Fossil fuels"),
    @(18, "Natural Gas", "OTH_ALL_NGS_001", "Fossil fuels"),
    @(19, "Natural Gas", "PWR_NGS_001", "Fossil fuels"),
    @(20, "LPG", "RES_CWH_LPG_001", "Fossil fuels"),
    @(21, "Oil products", "IND_ALL_OILPRD_001", "Fossil fuels"),
    @(22, "Oil products", "OTH_ALL_OILPRD_001", "Fossil fuels"),
    @(23, "Oil products", "TRA_AN_OILPRD_001", "Fossil fuels"),
    @(24, "Oil products", "PWR_OILPRD_001", "Fossil fuels"),
    @(25, "Diesel", "COM_CWH_OILPRD_001", "Fossil fuels"),
    @(26, "Diesel", "TRA_BUS_DSL_001", "Fossil fuels"),
    @(27, "Diesel", "TRA_CAR_DSL_001", "Fossil fuels"),
    @(28, "Diesel", "TRA_RLW_FREIGHT_DSL_001", "Fossil fuels"),
    @(29, "Diesel", "TRA_RLW_PSNG_DSL_001", "Fossil fuels"),
    @(30, "Diesel", "TRA_TRUCK_001", "Fossil fuels"),
    @(31, "Diesel", "PWR_OILPRD_002", "Fossil fuels"),
    @(32, "Diesel", "PWR_OILPRD_003", "Fossil fuels"),
    @(33, "Gasoline", "TRA_CAR_GSL_001", "Fossil fuels"),
    @(34, "Gasoline", "TRA_MCY_001", "Fossil fuels"),
    @(35, "Coal", "IND_ALL_COA_001", "Fossil fuels"),
    @(36, "Coal", "PWR_COA_001", "Fossil fuels")
)

foreach ($row in $target) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# Move the active selection, matching the author's final cursor position.
$ws.Range("C8").Select()
